$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Numeric cells (Duration / avg. fps) - plain numeric assignment
$numCells = @{
    "H3" = 30
    "H4" = 30
    "H5" = 30
    "H6" = 30
    "H7" = 30
    "H8" = 30
    "H9" = 30
    "H10" = 30
    "H11" = 30
    "H12" = 30
    "H13" = 30
    "H14" = 30
    "H15" = 30
    "H16" = 30
    "H17" = 30
    "H18" = 30
    "H19" = 30
    "H20" = 30
    "H21" = 30
    "H22" = 30
    "H23" = 30
    "H24" = 30
    "H25" = 30
    "H26" = 30
    "H27" = 30
    "H28" = 30
    "B29" = 28
    "H29" = 30
    "B30" = 29
    "H30" = 30
    "B31" = 30
    "H31" = 30
    "B32" = 31
    "H32" = 30
    "B33" = 32
    "B34" = 33
    "B35" = 34
    "B36" = 35
    "B37" = 36
    "B38" = 37
    "B39" = 38
    "B40" = 39
    "B41" = 40
    "B42" = 41
    "B43" = 42
    "B44" = 43
    "B45" = 44
    "B46" = 45
    "B47" = 46
    "B48" = 47
    "B49" = 48
    "B50" = 49
    "B51" = 50
    "H51" = 29.98
    "B52" = 51
    "H52" = 29.98
    "B53" = 52
    "H53" = 29.98
    "B54" = 53
    "H54" = 29.98
    "B55" = 54
    "H55" = 29.98
    "B56" = 55
    "H56" = 29.98
    "B57" = 56
    "H57" = 29.98
    "B58" = 57
    "H58" = 29.98
    "B59" = 58
    "H59" = 29.98
    "H87" = 30
    "H88" = 30
    "H89" = 30
    "B90" = 29
    "H90" = 30
    "B91" = 30
    "H91" = 30
    "B92" = 31
    "H92" = 30
    "B93" = 32
    "B94" = 33
    "B95" = 34
    "B96" = 35
    "B97" = 36
    "B98" = 37
    "B99" = 38
    "B100" = 39
    "H100" = 30
    "B101" = 40
    "H101" = 30
    "B102" = 41
    "H102" = 30
    "B103" = 42
    "H103" = 30
    "B104" = 43
    "H104" = 30
    "B105" = 44
    "H105" = 30
    "B106" = 45
    "H106" = 30
    "B107" = 46
    "H107" = 30
    "B108" = 47
    "B109" = 48
    "B110" = 49
}
foreach ($addr in $numCells.Keys) {
    $ws.Range($addr).Value = $numCells[$addr]
}

# Text-typed cells (MB/s / # frames / # Total Frame) - force text so shared-string type is preserved
$strCells = @{
    "I3" = "295.38"
    "I4" = "295.40"
    "I5" = "295.43"
    "I6" = "295.41"
    "I8" = "295.33"
    "I9" = "295.50"
    "I10" = "295.40"
    "I11" = "295.59"
    "I12" = "295.39"
    "I13" = "295.39"
    "I14" = "295.44"
    "I15" = "295.42"
    "I16" = "295.41"
    "I17" = "295.42"
    "I18" = "295.42"
    "I19" = "295.39"
    "I20" = "295.37"
    "I21" = "295.44"
    "I24" = "295.48"
    "I25" = "295.49"
    "I26" = "295.50"
    "I27" = "295.41"
    "I28" = "295.43"
    "I29" = "295.42"
    "I30" = "295.30"
    "I31" = "295.49"
    "I32" = "295.43"
    "I34" = "295.38"
    "I36" = "295.40"
    "I37" = "295.44"
    "I39" = "295.40"
    "I40" = "295.42"
    "I41" = "295.63"
    "I42" = "295.43"
    "I43" = "295.39"
    "I44" = "295.40"
    "I45" = "295.42"
    "I46" = "295.44"
    "I47" = "295.41"
    "I48" = "295.39"
    "I49" = "295.40"
    "I50" = "295.43"
    "C51" = "1499"
    "F51" = "29"
    "I51" = "295.42"
    "C52" = "1529"
    "I52" = "295.41"
    "C53" = "1559"
    "I53" = "295.38"
    "C54" = "1589"
    "I54" = "295.44"
    "C55" = "1619"
    "I55" = "295.42"
    "C56" = "1649"
    "I56" = "295.49"
    "C57" = "1679"
    "I57" = "295.41"
    "C58" = "1709"
    "I58" = "295.49"
    "C59" = "1739"
    "F60" = "31"
    "I60" = "295.37"
    "I61" = "295.44"
    "I63" = "295.43"
    "I64" = "295.39"
    "I67" = "295.42"
    "I68" = "295.41"
    "I69" = "295.59"
    "I70" = "295.45"
    "I71" = "295.37"
    "I72" = "295.41"
    "I74" = "295.41"
    "I75" = "295.40"
    "I76" = "295.41"
    "I77" = "295.40"
    "I78" = "295.43"
    "I79" = "295.42"
    "I80" = "295.39"
    "I81" = "295.41"
    "I82" = "295.41"
    "I83" = "295.39"
    "I84" = "295.49"
    "I85" = "295.51"
    "I86" = "295.45"
    "I87" = "295.39"
    "I88" = "295.44"
    "I90" = "295.36"
    "I91" = "295.47"
    "I92" = "295.42"
    "I93" = "295.39"
    "I96" = "295.43"
    "I97" = "295.37"
    "I98" = "295.45"
    "I99" = "295.40"
    "I100" = "295.43"
    "I101" = "295.56"
    "I102" = "295.45"
    "I103" = "295.41"
    "I105" = "295.39"
    "I106" = "295.31"
    "I107" = "295.50"
    "I108" = "295.42"
    "I109" = "295.41"
    "I110" = "295.42"
}
foreach ($addr in $strCells.Keys) {
    $rng = $ws.Range($addr)
    $rng.NumberFormat = "@"
    $rng.Value = $strCells[$addr]
    $rng.ClearFormats()
}
